# Stabilize upload flow, add structured OCR & payment analytics, fix redirects
#
# Rebuilds the payments sheet's header row and sample data row:
#   - new columns inserted for invoice_no, amount_detected, currency
#   - fraud/status/processing metadata appended (fraud_risk_score,
#     final_status, processed_utc)
#   - existing "signature_present" header shifted into place; the old
#     ink_fraction/amount/payable columns are retired
#   - sample row updated to match the new schema

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
# Copy the formatting of the existing bold/bordered header cell (A1) onto
# the three brand-new header cells so they pick up the same style (s="1")
# instead of minting a duplicate style entry.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

$ws.Range("A1").Value = "file"
$ws.Range("B1").Value = "invoice_no"
$ws.Range("C1").Value = "amount_detected"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "signature_present"
$ws.Range("F1").Value = "fraud_risk_score"
$ws.Range("G1").Value = "final_status"
$ws.Range("H1").Value = "processed_utc"

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "sample1.PNG"
$ws.Range("B2").Value = "IOCL-2024-INV-001"
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = "INR"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 0.08
$ws.Range("G2").Value = "APPROVED"
$ws.Range("H2").Value = "2025-12-17 08:36:40"
